# Applies the BanList.docx edit described by the commit:
#  1. Picture 1 (first inline image) resized: 496.35pt x 291.25pt -> 496.5pt x 291pt,
#     and a trailing run containing a single space is appended after it.
#  2. The "IP bans will work..." sentence loses its <w:proofErr/> gramStart/gramEnd
#     markers (and the three runs collapse into the two runs Word left behind).
#  3. Picture 3 (second inline image) loses its <w:lastRenderedPageBreak/> and is
#     resized: height 138.55pt -> 138.75pt.

$d = $word.ActiveDocument

function Get-ParagraphByMarker([string]$marker) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $rng = $p.Range
        if ($rng.Text -like "*$marker*") {
            return $rng
        }
    }
    return $null
}

# --- 1) Picture 1: resize + add trailing space run ---------------------
$picPara = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i).Range
    if ($cand.WordOpenXML -like "*Picture 1*") {
        $picPara = $cand
    }
}
$xml = $picPara.WordOpenXML
$xml = $xml -replace "width:496\.35pt;height:291\.25pt", "width:496.5pt;height:291pt"
$xml = $xml -replace "</w:pict></w:r></w:p>", "</w:pict></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r></w:p>"
$picPara.InsertXML($xml)

# --- 2) Drop the gramStart/gramEnd proofErr markers around "bans" ------
$bansPara = Get-ParagraphByMarker("you save the screen")
$xml2 = $bansPara.WordOpenXML
$bansPara.InsertXML($xml2)

# --- 3) Picture 3: drop lastRenderedPageBreak + resize height ----------
$picPara3 = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i).Range
    if ($cand.WordOpenXML -like "*Picture 3*") {
        $picPara3 = $cand
    }
}
$xml3 = $picPara3.WordOpenXML
$xml3 = $xml3 -replace "height:138\.55pt", "height:138.75pt"
$picPara3.InsertXML($xml3)
